$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values
$ws.Range("A2").Value = "carlos"
$ws.Range("B2").Value = "carro legal"
$ws.Range("C2").Value = "Tokio"
$ws.Range("D2").Value = 1200
$ws.Range("E2").Value = 10
$ws.Range("F2").NumberFormat = "@"
$ws.Range("F2").Value = "120.00"
$ws.Range("G2").Value = 1500
$ws.Range("H2").Value = "15/07/2004"
$ws.Range("I2").Value = "15/07/2006"
$ws.Range("J2").Value = $false
# K2 / L2 stay untouched (already empty string cells)

# Remove rows 3, 4, 5 entirely
$ws.Rows("3:5").Delete()
